$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 281.7688496666667
$ws.Range("H2").Value = 845.306549
$ws.Range("I2").Value = 0.4678027549763871
$ws.Range("J2").Value = 0.4678027549763872
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 68.00339
$ws.Range("N2").Value = 204.01017
$ws.Range("O2").Value = 0.6265962299909886
$ws.Range("P2").Value = 0.6265962299909885
$ws.Range("Q2").Value = 19161.2369737337
$ws.Range("R2").Value = 172451.1327636033
$ws.Range("S2").Value = 0.2931234426476023
$ws.Range("T2").Value = 0.2931234426476024
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 281.7688496666667
$ws.Range("H3").Value = 845.306549
$ws.Range("I3").Value = 0.4678027549763871
$ws.Range("J3").Value = 0.4678027549763872
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.991529999999999
$ws.Range("N3").Value = 26.97459
$ws.Range("O3").Value = 0.08284967558015671
$ws.Range("P3").Value = 0.08284967558015671
$ws.Range("Q3").Value = 2533.533064843323
$ws.Range("R3").Value = 22801.79758358991
$ws.Range("S3").Value = 0.03875730648529721
$ws.Range("T3").Value = 0.03875730648529722
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 281.7688496666667
$ws.Range("H4").Value = 845.306549
$ws.Range("I4").Value = 0.4678027549763871
$ws.Range("J4").Value = 0.4678027549763872
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.76843933333333
$ws.Range("N4").Value = 32.305318
$ws.Range("O4").Value = 0.09922245772090688
$ws.Range("P4").Value = 0.09922245772090688
$ws.Range("Q4").Value = 3034.21076365862
$ws.Range("R4").Value = 27307.89687292758
$ws.Range("S4").Value = 0.04641653907736833
$ws.Range("T4").Value = 0.04641653907736834
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 281.7688496666667
$ws.Range("H5").Value = 845.306549
$ws.Range("I5").Value = 0.4678027549763871
$ws.Range("J5").Value = 0.4678027549763872
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 20.764887
$ws.Range("N5").Value = 62.294661
$ws.Range("O5").Value = 0.1913316367079478
$ws.Range("P5").Value = 0.1913316367079478
$ws.Range("Q5").Value = 5850.89832344832
$ws.Range("R5").Value = 52658.08491103489
$ws.Range("S5").Value = 0.0895054667661192
$ws.Range("T5").Value = 0.08950546676611923
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 96.77942399999999
$ws.Range("H6").Value = 290.338272
$ws.Range("I6").Value = 0.160676672477411
$ws.Range("J6").Value = 0.1606766724774111
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 68.00339
$ws.Range("N6").Value = 204.01017
$ws.Range("O6").Value = 0.6265962299909886
$ws.Range("P6").Value = 0.6265962299909885
$ws.Range("Q6").Value = 6581.328914247359
$ws.Range("R6").Value = 59231.96022822623
$ws.Range("S6").Value = 0.1006793972218426
$ws.Range("T6").Value = 0.1006793972218426
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 96.77942399999999
$ws.Range("H7").Value = 290.338272
$ws.Range("I7").Value = 0.160676672477411
$ws.Range("J7").Value = 0.1606766724774111
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.991529999999999
$ws.Range("N7").Value = 26.97459
$ws.Range("O7").Value = 0.08284967558015671
$ws.Range("P7").Value = 0.08284967558015671
$ws.Range("Q7").Value = 870.1950942787198
$ws.Range("R7").Value = 7831.755848508478
$ws.Range("S7").Value = 0.0133120101880526
$ws.Range("T7").Value = 0.0133120101880526
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 96.77942399999999
$ws.Range("H8").Value = 290.338272
$ws.Range("I8").Value = 0.160676672477411
$ws.Range("J8").Value = 0.1606766724774111
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.76843933333333
$ws.Range("N8").Value = 32.305318
$ws.Range("O8").Value = 0.09922245772090688
$ws.Range("P8").Value = 0.09922245772090688
$ws.Range("Q8").Value = 1042.163356058944
$ws.Range("R8").Value = 9379.470204530495
$ws.Range("S8").Value = 0.01594273434162592
$ws.Range("T8").Value = 0.01594273434162592
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 96.77942399999999
$ws.Range("H9").Value = 290.338272
$ws.Range("I9").Value = 0.160676672477411
$ws.Range("J9").Value = 0.1606766724774111
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.764887
$ws.Range("N9").Value = 62.294661
$ws.Range("O9").Value = 0.1913316367079478
$ws.Range("P9").Value = 0.1913316367079478
$ws.Range("Q9").Value = 2009.613803285088
$ws.Range("R9").Value = 18086.52422956579
$ws.Range("S9").Value = 0.03074253072588992
$ws.Range("T9").Value = 0.03074253072588993
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 141.6283213333333
$ws.Range("H10").Value = 424.884964
$ws.Range("I10").Value = 0.2351364211508588
$ws.Range("J10").Value = 0.2351364211508588
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 68.00339
$ws.Range("N10").Value = 204.01017
$ws.Range("O10").Value = 0.6265962299909886
$ws.Range("P10").Value = 0.6265962299909885
$ws.Range("Q10").Value = 9631.205970675986
$ws.Range("R10").Value = 86680.85373608387
$ws.Range("S10").Value = 0.1473355950267015
$ws.Range("T10").Value = 0.1473355950267015
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 141.6283213333333
$ws.Range("H11").Value = 424.884964
$ws.Range("I11").Value = 0.2351364211508588
$ws.Range("J11").Value = 0.2351364211508588
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.991529999999999
$ws.Range("N11").Value = 26.97459
$ws.Range("O11").Value = 0.08284967558015671
$ws.Range("P11").Value = 0.08284967558015671
$ws.Range("Q11").Value = 1273.455300118306
$ws.Range("R11").Value = 11461.09770106476
$ws.Range("S11").Value = 0.01948097620942775
$ws.Range("T11").Value = 0.01948097620942775
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 141.6283213333333
$ws.Range("H12").Value = 424.884964
$ws.Range("I12").Value = 0.2351364211508588
$ws.Range("J12").Value = 0.2351364211508588
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.76843933333333
$ws.Range("N12").Value = 32.305318
$ws.Range("O12").Value = 0.09922245772090688
$ws.Range("P12").Value = 0.09922245772090688
$ws.Range("Q12").Value = 1525.115986159839
$ws.Range("R12").Value = 13726.04387543855
$ws.Range("S12").Value = 0.02333081360628644
$ws.Range("T12").Value = 0.02333081360628644
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 141.6283213333333
$ws.Range("H13").Value = 424.884964
$ws.Range("I13").Value = 0.2351364211508588
$ws.Range("J13").Value = 0.2351364211508588
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 20.764887
$ws.Range("N13").Value = 62.294661
$ws.Range("O13").Value = 0.1913316367079478
$ws.Range("P13").Value = 0.1913316367079478
$ws.Range("Q13").Value = 2940.896088486356
$ws.Range("R13").Value = 26468.0647963772
$ws.Range("S13").Value = 0.04498903630844313
$ws.Range("T13").Value = 0.04498903630844314
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 82.14745433333333
$ws.Range("H14").Value = 246.442363
$ws.Range("I14").Value = 0.1363841513953429
$ws.Range("J14").Value = 0.1363841513953429
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 68.00339
$ws.Range("N14").Value = 204.01017
$ws.Range("O14").Value = 0.6265962299909886
$ws.Range("P14").Value = 0.6265962299909885
$ws.Range("Q14").Value = 5586.305374536856
$ws.Range("R14").Value = 50276.74837083171
$ws.Range("S14").Value = 0.0854577950948421
$ws.Range("T14").Value = 0.0854577950948421
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 82.14745433333333
$ws.Range("H15").Value = 246.442363
$ws.Range("I15").Value = 0.1363841513953429
$ws.Range("J15").Value = 0.1363841513953429
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 8.991529999999999
$ws.Range("N15").Value = 26.97459
$ws.Range("O15").Value = 0.08284967558015671
$ws.Range("P15").Value = 0.08284967558015671
$ws.Range("Q15").Value = 738.6313000617965
$ws.Range("R15").Value = 6647.68170055617
$ws.Range("S15").Value = 0.01129938269737914
$ws.Range("T15").Value = 0.01129938269737914
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 82.14745433333333
$ws.Range("H16").Value = 246.442363
$ws.Range("I16").Value = 0.1363841513953429
$ws.Range("J16").Value = 0.1363841513953429
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.76843933333333
$ws.Range("N16").Value = 32.305318
$ws.Range("O16").Value = 0.09922245772090688
$ws.Range("P16").Value = 0.09922245772090688
$ws.Range("Q16").Value = 884.5998783762703
$ws.Range("R16").Value = 7961.398905386434
$ws.Range("S16").Value = 0.01353237069562617
$ws.Range("T16").Value = 0.01353237069562618
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 82.14745433333333
$ws.Range("H17").Value = 246.442363
$ws.Range("I17").Value = 0.1363841513953429
$ws.Range("J17").Value = 0.1363841513953429
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 20.764887
$ws.Range("N17").Value = 62.294661
$ws.Range("O17").Value = 0.1913316367079478
$ws.Range("P17").Value = 0.1913316367079478
$ws.Range("Q17").Value = 1705.782606569327
$ws.Range("R17").Value = 15352.04345912394
$ws.Range("S17").Value = 0.0260946029074955
$ws.Range("T17").Value = 0.0260946029074955

Write-Host "Applied 224 cell updates"
